# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the other columns, and fill in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: text + same formatting as the neighboring "sum" header (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row "Save" values.
$values = @(0, 0, 0, 0, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
